$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (and values with multiple dots) are preserved exactly as literal text,
# matching the inline-string cell type used in the source workbook.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.141.26'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '2.566.60'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '316.75'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').Value = '96.72'
$ws.Range('E6').Value = '  +1.75%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('D10').Value = '35.69'
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '7.46'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '2.964.04'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('E14').Value = '  -4.62%  '
$ws.Range('D15').Value = '2.571.64'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '15.04'
$ws.Range('E16').Value = '  -2.96%  '
$ws.Range('D17').Value = '0.847'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '43.167.16'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = '6.87'
$ws.Range('E19').Value = '  +4.58%  '
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').Value = '69.57'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '253.15'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').Value = '2.96'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  +2.48%  '
$ws.Range('D26').Value = '26.83'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  +1.73%  '
$ws.Range('D29').Value = '40.15'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').Value = '5.83'
$ws.Range('E31').Value = '  -3.93%  '
$ws.Range('D32').Value = '154.24'
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('E33').Value = '  +3.97%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '2.14'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0809'
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('D36').Value = '2.70'
$ws.Range('E36').Value = '  +2.86%  '
$ws.Range('D37').Value = '19.05'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('E39').Value = '  +5.11%  '
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').Value = '22.51'
$ws.Range('E41').Value = '  -5.07%  '
$ws.Range('D42').Value = '3.89'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '3.28'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = '1.994.80'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '9.03'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').Value = '84.01'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('D49').Value = '2.817.67'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').Value = '74.19'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').Value = '104.13'

# Restore the default (unformatted) style so the saved XML matches the
# original "no explicit style" cells once more.
$dataRange.Style = "Normal"
